$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 1021.934916333333
$ws.Range("N2").Value = 3065.804749
$ws.Range("O2").Value = 0.8026347959846111
$ws.Range("P2").Value = 0.802634795984611
$ws.Range("Q2").Value = 205858.3065682682
$ws.Range("R2").Value = 1852724.759114414
$ws.Range("S2").Value = 0.3879535473474867
$ws.Range("T2").Value = 0.3879535473474867
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.04931810976893385
$ws.Range("P3").Value = 0.04931810976893384
$ws.Range("Q3").Value = 12649.01872055812
$ws.Range("R3").Value = 113841.168485023
$ws.Range("S3").Value = 0.02383790950635221
$ws.Range("T3").Value = 0.0238379095063522
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("M4").Value = 187.139577
$ws.Range("N4").Value = 561.418731
$ws.Range("O4").Value = 0.1469807262726385
$ws.Range("P4").Value = 0.1469807262726385
$ws.Range("Q4").Value = 37697.34823362886
$ws.Range("R4").Value = 339276.1341026597
$ws.Range("S4").Value = 0.07104313746980055
$ws.Range("T4").Value = 0.07104313746980055
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("M5").Value = 1.357726666666667
$ws.Range("N5").Value = 4.073180000000001
$ws.Range("O5").Value = 0.001066367973816652
$ws.Range("P5").Value = 0.001066367973816652
$ws.Range("Q5").Value = 273.5001103450045
$ws.Range("R5").Value = 2461.500993105041
$ws.Range("S5").Value = 0.0005154289849998652
$ws.Range("T5").Value = 0.0005154289849998651
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 1021.934916333333
$ws.Range("N6").Value = 3065.804749
$ws.Range("O6").Value = 0.8026347959846111
$ws.Range("P6").Value = 0.802634795984611
$ws.Range("Q6").Value = 66852.28908737715
$ws.Range("R6").Value = 601670.6017863944
$ws.Range("S6").Value = 0.1259875451814556
$ws.Range("T6").Value = 0.1259875451814557
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("O7").Value = 0.04931810976893385
$ws.Range("P7").Value = 0.04931810976893384
$ws.Range("S7").Value = 0.007741338419243755
$ws.Range("T7").Value = 0.007741338419243755
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("M8").Value = 187.139577
$ws.Range("N8").Value = 561.418731
$ws.Range("O8").Value = 0.1469807262726385
$ws.Range("P8").Value = 0.1469807262726385
$ws.Range("Q8").Value = 12242.17795217475
$ws.Range("R8").Value = 110179.6015695727
$ws.Range("S8").Value = 0.02307119126247332
$ws.Range("T8").Value = 0.02307119126247332
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("M9").Value = 1.357726666666667
$ws.Range("N9").Value = 4.073180000000001
$ws.Range("O9").Value = 0.001066367973816652
$ws.Range("P9").Value = 0.001066367973816652
$ws.Range("Q9").Value = 88.81890047099114
$ws.Range("R9").Value = 799.3701042389201
$ws.Range("S9").Value = 0.0001673850722064368
$ws.Range("T9").Value = 0.0001673850722064367
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 1021.934916333333
$ws.Range("N10").Value = 3065.804749
$ws.Range("O10").Value = 0.8026347959846111
$ws.Range("P10").Value = 0.802634795984611
$ws.Range("Q10").Value = 61760.47452159828
$ws.Range("R10").Value = 555844.2706943845
$ws.Range("S10").Value = 0.1163916850184147
$ws.Range("T10").Value = 0.1163916850184147
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("O11").Value = 0.04931810976893385
$ws.Range("P11").Value = 0.04931810976893384
$ws.Range("Q11").Value = 3794.888879818793
$ws.Range("R11").Value = 34153.99991836914
$ws.Range("S11").Value = 0.007151718224335998
$ws.Range("T11").Value = 0.007151718224335997
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("M12").Value = 187.139577
$ws.Range("N12").Value = 561.418731
$ws.Range("O12").Value = 0.1469807262726385
$ws.Range("P12").Value = 0.1469807262726385
$ws.Range("Q12").Value = 11309.75064318212
$ws.Range("R12").Value = 101787.755788639
$ws.Range("S12").Value = 0.02131397054013438
$ws.Range("T12").Value = 0.02131397054013438
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("M13").Value = 1.357726666666667
$ws.Range("N13").Value = 4.073180000000001
$ws.Range("O13").Value = 0.001066367973816652
$ws.Range("P13").Value = 0.001066367973816652
$ws.Range("Q13").Value = 82.05399567403558
$ws.Range("R13").Value = 738.4859610663202
$ws.Range("S13").Value = 0.0001546361632253139
$ws.Range("T13").Value = 0.0001546361632253139
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 1021.934916333333
$ws.Range("N14").Value = 3065.804749
$ws.Range("O14").Value = 0.8026347959846111
$ws.Range("P14").Value = 0.802634795984611
$ws.Range("Q14").Value = 91427.96083784135
$ws.Range("R14").Value = 822851.6475405722
$ws.Range("S14").Value = 0.172302018437254
$ws.Range("T14").Value = 0.172302018437254
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("O15").Value = 0.04931810976893385
$ws.Range("P15").Value = 0.04931810976893384
$ws.Range("Q15").Value = 5617.815513491515
$ws.Range("R15").Value = 50560.33962142363
$ws.Range("S15").Value = 0.01058714361900189
$ws.Range("T15").Value = 0.01058714361900189
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("M16").Value = 187.139577
$ws.Range("N16").Value = 561.418731
$ws.Range("O16").Value = 0.1469807262726385
$ws.Range("P16").Value = 0.1469807262726385
$ws.Range("Q16").Value = 16742.54362357588
$ws.Range("R16").Value = 150682.8926121829
$ws.Range("S16").Value = 0.03155242700023027
$ws.Range("T16").Value = 0.03155242700023026
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("M17").Value = 1.357726666666667
$ws.Range("N17").Value = 4.073180000000001
$ws.Range("O17").Value = 0.001066367973816652
$ws.Range("P17").Value = 0.001066367973816652
$ws.Range("Q17").Value = 121.4697516686112
$ws.Range("R17").Value = 1093.2277650175
$ws.Range("S17").Value = 0.0002289177533850362
$ws.Range("T17").Value = 0.0002289177533850361
